$d = $word.ActiveDocument

# --- Row "Performance Requirements" -------------------------------------
# Implementation Constraint cell: swap the "developer is new to the
# languages" text for the laptop-scaling limitation.
$ok = $d.Content.Find.Execute(
    "Developer is new to learning the languages required to write the application. This will impact the speed of production.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The app is being run on a laptop which doesn" + [char]0x2019 + "t allow the app to be run at a much larger scale e.g taking in data of thousands of users each with thousands of transactions.",
    2)
if (-not $ok) { throw "Could not find the 'Developer is new to learning...' text" }

# Solution cell: swap the "access resources" text for the cloud
# deployment solution.
$ok = $d.Content.Find.Execute(
    "Access as many resources as possible to aid the developer in their knowledge and understanding.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deploying on a cloud based server in order to be able to scale up as it can handle greater volumes of traffic.",
    2)
if (-not $ok) { throw "Could not find the 'Access as many resources...' text" }

# --- Row "Time Limitations" ---------------------------------------------
# Implementation Constraint cell: trim the trailing "developer is slow"
# sentence, keeping the first sentence intact.
$ok = $d.Content.Find.Execute(
    "Less than a week to complete the project including design and implementation. This is a problem as the developer is slow.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Less than a week to complete the project including design and implementation.",
    2)
if (-not $ok) { throw "Could not find the 'Less than a week...' text" }

# Solution cell: insert a new paragraph ahead of the existing "Discussing
# an extension..." paragraph with the MVP-prioritisation mitigation.
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(7, 3)
$cell.Range.InsertParagraphBefore() | Out-Null

$cell = $tbl.Cell(7, 3)
$newPara = $cell.Range.Paragraphs.Item(1)
$newPara.Range.Text = "Prioritising focus on the MVP of the brief and choosing not to include the extensions."
